$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 (shifts rows 28.. down by one, including merged cells)
$ws.Rows("28:28").Insert()

# Copy formatting only from the row just below (the shifted former row 28) so the
# new row matches the existing item-row styling exactly.
$ws.Range("A29:Q29").Copy()
$ws.Range("A28:Q28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row with the new product's data
$ws.Range("A28").Value = 22
$ws.Range("C28").Value = "شيلز حريمي مزيل عرق "
$ws.Range("H28").Value = "6:0"
$ws.Range("L28").Value = "0"
$ws.Range("N28").Value = "75.00"
$ws.Range("P28").Value = "75.0000"
$ws.Range("Q28").Value = "1:0"

# The "#" column (A) is a positional row index (row-6), independent of which
# product occupies the row, so after the insert it must be re-stamped for
# every item row that shifted down (it is not part of the content shift).
$ws.Range("A29").Value = 23
$ws.Range("A30").Value = 24
$ws.Range("A31").Value = 25
$ws.Range("A32").Value = 26
$ws.Range("A33").Value = 27

# Update the running total on the totals row (now shifted down to row 34)
$ws.Range("P34").Value = 1571.12

# Update the generated timestamp shown in the footer (now shifted down to row 35)
$ws.Range("A35").Value = "Thursday, 11 September, 2025 1:58 PM"
